$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename distance condition codes and the "S30" size code across every cell
# that uses them (shared strings: condition names, filenames, the Distance
# column values, and the Size column values).
# Order matters: do the 3-digit distance codes first (they are unambiguous
# substrings), then the size code.
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
